# Applies the "Add commands for resetting git credentials" edit:
#  - Fix capitalization in the existing "Stash uncommitted changes" and
#    "Add a remote" command cells.
#  - Insert a new row for "Reset git credentials" with Mac/Windows
#    instructions after "Delete a branch (local+remote)".
#  - Update the view selection/scroll position to match (B9 selected,
#    top row scrolled to A7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 9 for "Reset git credentials" ----------------------

$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "Reset git credentials"
$ws.Range("B9").Value = "(For Mac)`ngit config --global credential.helper osxkeychain`n(Credentials will be asked on next pull/push)`n(For Windows, reset from Windows Credentials Manager)"

# Match formatting of the rows above: wrap text in column B, and a row
# height tall enough for the 4-line instructions (same as other 4-line
# rows, e.g. row 2/3/7 which use ht="58").
$ws.Range("B9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 58

# --- Fix capitalization of existing cells -----------------------------

$ws.Range("B6").Value = "git remote add origin https://github.com/<username>/<repo-name>.git`n(Make a commit)`ngit push --set-upstream origin master"

$ws.Range("B4").Value = "git stash`ngit stash push -m ""<stash-name>""`n(Switch branches)`ngit stash apply`ngit stash list`ngit stash apply stash@{1}"

# --- Update the view to match the diff ---------------------------------

$ws.Range("B9").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
